# Add the new "homework" worksheet after the last existing sheet.
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "homework"

# NOTE: shared-string indices are assigned in the order string VALUES are
# first written (not in sheet/row order), so the label writes below are
# deliberately sequenced to reproduce the target shared-string table:
# 61 숙제 5, 62 N, 63 n, 64 y_s, 65 V(y_s), 66 s^2, 67 ((N-n)/N)*s^2/n, 68 신뢰구간

$ws.Range("B2").Value = "숙제 5"
$ws.Range("B3").Value = "N"
$ws.Range("B4").Value = "n"
$ws.Range("B5").Value = "y_s"
$ws.Range("B8").Value = "V(y_s)"
$ws.Range("B6").Value = "s^2"
$ws.Range("D8").Value = "((N-n)/N)*s^2/n"
$ws.Range("B9").Value = "신뢰구간"

# ---- Block 1 numeric / formula content (rows 2-9) ----
$ws.Range("C3").Value = 1800
$ws.Cells.Item(3, 5).Value = 12
$ws.Cells.Item(3, 6).Value = 11.97
$ws.Cells.Item(3, 7).Value = 12.01
$ws.Cells.Item(3, 8).Value = 12.03
$ws.Cells.Item(3, 9).Value = 12.01
$ws.Cells.Item(3, 10).Value = 11.8

$ws.Range("C4").Formula = "=COUNT(E3:J8)"
$ws.Cells.Item(4, 5).Value = 11.91
$ws.Cells.Item(4, 6).Value = 11.98
$ws.Cells.Item(4, 7).Value = 12.03
$ws.Cells.Item(4, 8).Value = 11.98
$ws.Cells.Item(4, 9).Value = 12
$ws.Cells.Item(4, 10).Value = 11.83

$ws.Range("C5").Formula = "=SUM(E3:J8)/C4"
$ws.Cells.Item(5, 5).Value = 11.87
$ws.Cells.Item(5, 6).Value = 12.01
$ws.Cells.Item(5, 7).Value = 11.98
$ws.Cells.Item(5, 8).Value = 11.87
$ws.Cells.Item(5, 9).Value = 11.9
$ws.Cells.Item(5, 10).Value = 11.88

$ws.Range("C6").Formula = "=VAR.S(E3:J8)"
$ws.Cells.Item(6, 5).Value = 12.05
$ws.Cells.Item(6, 6).Value = 11.87
$ws.Cells.Item(6, 7).Value = 11.91
$ws.Cells.Item(6, 8).Value = 11.93
$ws.Cells.Item(6, 9).Value = 11.94
$ws.Cells.Item(6, 10).Value = 11.89

$ws.Cells.Item(7, 5).Value = 11.75
$ws.Cells.Item(7, 6).Value = 11.93
$ws.Cells.Item(7, 7).Value = 11.95
$ws.Cells.Item(7, 8).Value = 11.97
$ws.Cells.Item(7, 9).Value = 11.93
$ws.Cells.Item(7, 10).Value = 12.05

$ws.Range("C8").Formula = "=((C3-C4)/C3)*C6/C4"
$ws.Cells.Item(8, 5).Value = 11.85
$ws.Cells.Item(8, 6).Value = 11.98
$ws.Cells.Item(8, 7).Value = 11.87
$ws.Cells.Item(8, 8).Value = 12.05
$ws.Cells.Item(8, 9).Value = 12.02
$ws.Cells.Item(8, 10).Value = 12.04

$ws.Range("C9").Formula = "=C5-2*SQRT(C8)"
$ws.Range("D9").Formula = "=C5+2*SQRT(C8)"

# ---- Block 2 (rows 12-16) ----
$ws.Range("B12").Value = "N"
$ws.Range("C12").Formula = "=C13*100"
$ws.Cells.Item(12, 5).Value = 120
$ws.Cells.Item(12, 6).Value = 119.7
$ws.Cells.Item(12, 7).Value = 120.1
$ws.Cells.Item(12, 8).Value = 120.3

$ws.Range("B13").Value = "n"
$ws.Range("C13").Formula = "=COUNT(E12:H16)"
$ws.Cells.Item(13, 5).Value = 119.1
$ws.Cells.Item(13, 6).Value = 119.8
$ws.Cells.Item(13, 7).Value = 120.3
$ws.Cells.Item(13, 8).Value = 119.8

$ws.Range("B14").Value = "y_s"
$ws.Range("C14").Formula = "=SUM(E12:H16)/C13"
$ws.Cells.Item(14, 5).Value = 118.7
$ws.Cells.Item(14, 6).Value = 120.1
$ws.Cells.Item(14, 7).Value = 119.8
$ws.Cells.Item(14, 8).Value = 118.7

$ws.Range("B15").Value = "V(y_s)"
$ws.Range("C15").Formula = "=VAR.S(E12:H16)"
$ws.Cells.Item(15, 5).Value = 120.5
$ws.Cells.Item(15, 6).Value = 118.7
$ws.Cells.Item(15, 7).Value = 119.1
$ws.Cells.Item(15, 8).Value = 119.3

$ws.Cells.Item(16, 5).Value = 117.5
$ws.Cells.Item(16, 6).Value = 119.3
$ws.Cells.Item(16, 7).Value = 119.5
$ws.Cells.Item(16, 8).Value = 119.7

# ---- Sheet-level view / selection tweaks to match target file ----
$sheet1 = $wb.Worksheets.Item(1)
[void]$sheet1.Activate()
[void]$sheet1.Range("H14").Select()

$sheet2 = $wb.Worksheets.Item(2)
[void]$sheet2.Activate()
[void]$sheet2.Range("H13").Select()

[void]$ws.Activate()
[void]$ws.Range("D10").Select()

Write-Host "done"
